# Workhours.xlsx update: "Improved car visuals and camera."
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rewrite the task description that used to read "Importing rigged 3D
#    character" (row 14, the continuation of the second table) and give the
#    now-open "To" cell next to it a "~12:30" text marker.
# ---------------------------------------------------------------------------
$ws.Range("D14").Value = "Importing low poly asset pack for character and car visuals."
$ws.Range("C14").NumberFormat = "h:mm"
$ws.Range("C14").Value = "~12:30"

# ---------------------------------------------------------------------------
# 2) Append the new rows (15-21) describing the rest of the day.
# ---------------------------------------------------------------------------
$ws.Range("B15").NumberFormat = "h:mm"
$ws.Range("B15").Value = "~12:30"
$ws.Range("C15").NumberFormat = "h:mm"
$ws.Range("C15").Value = 0.55208333333333337
$ws.Range("D15").Value = "Modelling missing or unmovable props (steering wheel, handbrake, pedals, gearshift)"

$ws.Range("B16").NumberFormat = "h:mm"
$ws.Range("B16").Value = 0.55208333333333337
$ws.Range("C16").NumberFormat = "h:mm"
$ws.Range("C16").Value = 0.57638888888888895
$ws.Range("D16").Value = "Importing the new assets, fixing rotation and scale problems. There are some problems with the pedals normals."

$ws.Range("B17").NumberFormat = "h:mm"
$ws.Range("B17").Value = 0.57638888888888895
$ws.Range("C17").NumberFormat = "h:mm"
$ws.Range("C17").Value = 0.57638888888888895
$ws.Range("D17").Value = "Realized, that i was editing different files, then what I exported to Unity, cause I made a backup to different folders…"

$ws.Range("B18").NumberFormat = "h:mm"
$ws.Range("B18").Value = 0.57638888888888895
$ws.Range("C18").NumberFormat = "h:mm"
$ws.Range("C18").Value = 0.59027777777777779
$ws.Range("D18").Value = "Editing the real files… Still, there is some bug with the import, the rotation is off, but it does not really matter."

$ws.Range("B19").NumberFormat = "h:mm"
$ws.Range("B19").Value = 0.59027777777777779
$ws.Range("C19").NumberFormat = "h:mm"
$ws.Range("C19").Value = 0.60763888888888895
$ws.Range("D19").Value = "Finished the the car visuals."

$ws.Range("B20").NumberFormat = "h:mm"
$ws.Range("B20").Value = 0.61111111111111105
$ws.Range("C20").NumberFormat = "h:mm"
$ws.Range("C20").Value = 0.67361111111111116
$ws.Range("D20").Value = "Importing and stripping down my camera controller from another project. Added different camera positions with limited rotation."

$ws.Range("B21").NumberFormat = "h:mm"
$ws.Range("B21").Value = 0.67361111111111116
$ws.Range("C21").NumberFormat = "h:mm"
$ws.Range("C21").Value = 0.67708333333333337
$ws.Range("D21").Value = "Added headlights."

# ---------------------------------------------------------------------------
# 3) Word-wrap the whole "Task description" column (including the two table
#    headers), and give it a bit more width to stay readable.
# ---------------------------------------------------------------------------
$ws.Range("D2").WrapText = $true
$ws.Range("D3:D9").WrapText = $true
$ws.Range("D11").WrapText = $true
$ws.Range("D12:D21").WrapText = $true
$ws.Range("D16").VerticalAlignment = -4108

$ws.Columns.Item(4).ColumnWidth = 57.5

# ---------------------------------------------------------------------------
# 4) Rows whose wrapped text spills onto more than one line need an explicit
#    (taller) row height.
# ---------------------------------------------------------------------------
$ws.Rows.Item(15).RowHeight = 30
$ws.Rows.Item(16).RowHeight = 30.75
$ws.Rows.Item(17).RowHeight = 30
$ws.Rows.Item(18).RowHeight = 30
$ws.Rows.Item(20).RowHeight = 45

# ---------------------------------------------------------------------------
# 5) View state: narrower zoom now that there is more content, and the
#    selection the author left the sheet on.
# ---------------------------------------------------------------------------
$ws.Select()
$ws.Range("G19").Select()
$excel.ActiveWindow.Zoom = 145
